$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("HEATING")
$ws2 = $wb.Worksheets.Item("HOT_WATER")

# Row 38 on HEATING
$ws1.Range("A38").Value = "condensing boiler + solar thermal collector"
$ws1.Range("B38").Value = "SUPPLY_HEATING_AS36"
$ws1.Range("C38").Value = "BO1"
$ws1.Range("D38").Value = "SC1"
$ws1.Range("E38").Value = "-"
$ws1.Range("F38").Value = "NATURALGAS"
$ws1.Range("G38").Value = "BUILDING"
$ws1.Range("H38").Value = 1
$ws1.Range("I38").Formula = "=142+600"
$ws1.Range("J38").Value = 20
$ws1.Range("K38").Value = 3
$ws1.Range("L38").Value = 5
$ws1.Range("M38").Value = "KEA Technikatalog"

# Row 39 on HEATING
$ws1.Range("A39").Value = "solar thermal collector"
$ws1.Range("B39").Value = "SUPPLY_HEATING_AS37"
$ws1.Range("C39").Value = "SC1"
$ws1.Range("D39").Value = "-"
$ws1.Range("E39").Value = "-"
$ws1.Range("F39").Value = "SOLAR"
$ws1.Range("G39").Value = "BUILDING"
$ws1.Range("H39").Value = 1
$ws1.Range("I39").Value = 600
$ws1.Range("J39").Value = 20
$ws1.Range("K39").Value = 3
$ws1.Range("L39").Value = 5
$ws1.Range("M39").Value = "KEA Technikatalog"

# Row 38 on HOT_WATER
$ws2.Range("A38").Value = "Flatplate collector"
$ws2.Range("B38").Value = "SUPPLY_HOTWATER_AS36"
$ws2.Range("C38").Value = "SOLAR"
$ws2.Range("D38").Value = "BUILDING"
$ws2.Range("E38").Value = 1
$ws2.Range("F38").Value = 600
$ws2.Range("G38").Value = 20
$ws2.Range("H38").Value = 3
$ws2.Range("I38").Value = 5
$ws2.Range("J38").Value = "KEA Technikatalog"
